{"js": "// Commit: \"Put first deliverables in their own folder, moved GDD to Design\n// folder, created a \"Title Ideas\" doc and added my ideas to it\"\n//\n// In the document body this shows up as the list item that used to just say\n// \"Uploaded deliverables\" growing into \"Uploaded deliverables to Spaces\"\n// (the trailing single-space run that used to sit right after the\n// \"_GoBack\" bookmark is dropped once the new text takes its place).\n\nconst body = context.document.body;\n\n// Locate the paragraph that holds \"Uploaded deliverables\".\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\nconst targetParagraph = body.paragraphs.items.find(\n  (p) => p.text.indexOf(\"Uploaded deliverables\") !== -1\n);\n\nif (!targetParagraph) {\n  throw new Error('Could not find the \"Uploaded deliverables\" paragraph.');\n}\n\n// Find the exact run of text \"Uploaded deliverables\" inside that paragraph\n// and append \" to Spaces\" right after it.\nconst found = targetParagraph.search(\"Uploaded deliverables\", {\n  matchCase: true,\n});\nfound.load(\"items\");\nawait context.sync();\n\nconst match = found.items[0];\nconst inserted = match.insertText(\" to Spaces\", \"After\");\nawait context.sync();\n\n// The paragraph originally ended with: ... <bookmark _GoBack/> \" \" (a lone\n// trailing space run after the bookmark). That trailing space is no longer\n// wanted now that real text follows the bookmark, so remove everything\n// between the end of what we just inserted and the end of the paragraph.\nconst afterInserted = inserted.getRange(\"After\");\nconst paragraphEnd = targetParagraph.getRange(\"End\");\nconst trailing = afterInserted.expandTo(paragraphEnd);\ntrailing.load(\"text\");\nawait context.sync();\n\nif (trailing.text.length > 0) {\n  trailing.delete();\n  await context.sync();\n}\n", "ps1": "# Commit: \"Put first deliverables in their own folder, moved GDD to Design\n# folder, created a \"Title Ideas\" doc and added my ideas to it\"\n#\n# The list item that used to read just \"Uploaded deliverables\" grows into\n# \"Uploaded deliverables to Spaces\" (the lone trailing space run that used\n# to sit right after the \"_GoBack\" bookmark is dropped once real text\n# follows the bookmark).\n\n$d = $word.ActiveDocument\n\n# Find the exact text \"Uploaded deliverables\" and append \" to Spaces\"\n# right after it.\n$found = $d.Content\n$ok = $found.Find.Execute(\"Uploaded deliverables\")\nif (-not $ok) {\n    throw 'Could not find \"Uploaded deliverables\" in the document.'\n}\n$found.Collapse(0)\n$found.InsertAfter(\" to Spaces\")\n\n# The paragraph originally ended with: ... <bookmark _GoBack/> \" \" (a lone\n# trailing space run after the bookmark). Now that real text follows the\n# bookmark, drop that leftover trailing space - i.e. everything between the\n# end of the bookmark and the end of the paragraph.\n$bm = $d.Bookmarks(\"_GoBack\")\n$para = $bm.Range.Paragraphs(1)\n$paragraphEnd = $para.Range.End\nif ($paragraphEnd - 1 -gt $bm.End) {\n    $trailing = $d.Range($bm.End, $paragraphEnd - 1)\n    if ($trailing.Text -ne $null -and $trailing.Text.Length -gt 0) {\n        $trailing.Delete()\n    }\n}\n"}
